$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Executive Summary")
$chart = $ws.ChartObjects(1).Chart
$ser = $chart.SeriesCollection(1)
$ser.HasDataLabels = $true
Write-Host "done"
